$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.601.11"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "2.490.27"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.24"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.55"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("D9").Value = "2.490.63"
$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("E11").Value = "  +1.79%  "

$ws.Range("E12").Value = "  -1.27%  "

$ws.Range("E13").Value = "  -2.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.28"
$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("D15").Value = "2.953.28"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").Value = "67.762.83"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "2.483.84"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.80"
$ws.Range("E19").Value = "  +3.70%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.01"
$ws.Range("E21").Value = "  +3.11%  "

$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("E24").Value = "  +1.85%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").Value = "  -4.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "2.622.74"
$ws.Range("E29").Value = "  -2.14%  "

$ws.Range("D30").Value = "0.0₃0957"
$ws.Range("E30").Value = "  -3.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "531.79"
$ws.Range("E32").Value = "  -1.04%  "

$ws.Range("E33").Value = "  -4.35%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("E36").Value = "  -3.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.74"
$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.42"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.69"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("E40").Value = "  +1.05%  "

$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.69"
$ws.Range("E46").Value = "  -3.81%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0275"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.67"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.546"
$ws.Range("E49").Value = "  -3.17%  "

$ws.Range("E50").Value = "  -1.76%  "

$ws.Range("E51").Value = "  -1.74%  "
